$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Yearly")

# Update sun profile (SOLAR column C) values for rows 5 and 6
$ws.Range("C5").Value = 1100
$ws.Range("C6").Value = 1200

# Update grid connection selection (active cell moved from D12 to E6)
$ws.Range("E6").Select()
